# Update "Pagos" (col F) and "Inscrições homologadas" (col H) values
# on the "Inscricoes" worksheet, per the data refresh described in the diff.
# Column G ("Isenções deferidas") is unchanged; H = F + G for every row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$updates = @(
    @{ Row = 2; F = 58; H = 70 },
    @{ Row = 4; F = 11; H = 25 },
    @{ Row = 7; F = 10; H = 11 },
    @{ Row = 8; F = 12; H = 13 },
    @{ Row = 9; F = 21; H = 30 },
    @{ Row = 15; F = 122; H = 163 },
    @{ Row = 17; F = 87; H = 119 },
    @{ Row = 18; F = 81; H = 118 },
    @{ Row = 19; F = 55; H = 68 },
    @{ Row = 24; F = 23; H = 27 },
    @{ Row = 26; F = 28; H = 38 },
    @{ Row = 28; F = 21; H = 23 },
    @{ Row = 29; F = 14; H = 17 },
    @{ Row = 31; F = 2; H = 3 },
    @{ Row = 32; F = 9; H = 18 },
    @{ Row = 33; F = 18; H = 30 },
    @{ Row = 34; F = 17; H = 20 },
    @{ Row = 36; F = 75; H = 107 },
    @{ Row = 37; F = 46; H = 58 },
    @{ Row = 38; F = 30; H = 50 },
    @{ Row = 40; F = 23; H = 25 },
    @{ Row = 41; F = 30; H = 41 },
    @{ Row = 42; F = 34; H = 43 },
    @{ Row = 43; F = 26; H = 29 },
    @{ Row = 45; F = 17; H = 24 },
    @{ Row = 46; F = 16; H = 25 },
    @{ Row = 47; F = 55; H = 65 },
    @{ Row = 48; F = 37; H = 43 },
    @{ Row = 49; F = 54; H = 71 },
    @{ Row = 50; F = 14; H = 23 },
    @{ Row = 52; F = 8; H = 8 },
    @{ Row = 55; F = 5; H = 8 },
    @{ Row = 56; F = 6; H = 8 },
    @{ Row = 57; F = 10; H = 14 },
    @{ Row = 59; F = 8; H = 12 },
    @{ Row = 61; F = 22; H = 32 },
    @{ Row = 62; F = 21; H = 35 },
    @{ Row = 63; F = 23; H = 31 },
    @{ Row = 64; F = 25; H = 30 },
    @{ Row = 65; F = 21; H = 34 },
    @{ Row = 66; F = 31; H = 39 },
    @{ Row = 67; F = 30; H = 38 },
    @{ Row = 69; F = 11; H = 14 },
    @{ Row = 70; F = 33; H = 46 },
    @{ Row = 71; F = 26; H = 36 },
    @{ Row = 72; F = 35; H = 46 },
    @{ Row = 73; F = 23; H = 35 },
    @{ Row = 74; F = 10; H = 14 },
    @{ Row = 75; F = 13; H = 18 },
    @{ Row = 76; F = 26; H = 43 },
    @{ Row = 77; F = 36; H = 53 },
    @{ Row = 78; F = 26; H = 47 },
    @{ Row = 79; F = 35; H = 46 },
    @{ Row = 80; F = 21; H = 33 },
    @{ Row = 81; F = 17; H = 22 },
    @{ Row = 83; F = 4; H = 11 },
    @{ Row = 87; F = 11; H = 18 },
    @{ Row = 88; F = 25; H = 33 },
    @{ Row = 89; F = 28; H = 34 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 6).Value = $u.F
    $ws.Cells.Item($u.Row, 8).Value = $u.H
}

$wb.Save()
